$wb = $excel.ActiveWorkbook

# Add the new "CreateEvent" sheet at the very end of the workbook (after the
# current last sheet, "AccessionRegister"), mirroring the diff which appends
# it as the 43rd sheet and makes it the active tab.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "CreateEvent"

# Header row
$ws.Range("A1").Value = "EventName"
$ws.Range("B1").Value = "Description"
$ws.Range("C1").Value = "Incharge"

# Data row
$ws.Range("A2").Value = "Test"
$ws.Range("B2").Value = "Testing"
$ws.Range("C2").Value = "Test"

# Column widths matching the authored sheet
$ws.Columns.Item(1).ColumnWidth = 18.5703125
$ws.Columns.Item(2).ColumnWidth = 32
$ws.Columns.Item(3).ColumnWidth = 21.28515625

# Select the same cell the authored sheet ends up with
[void]$ws.Range("C2").Select()
